$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.417.85"
$ws.Range("E2").Value = "  -3.59%  "
$ws.Range("D3").Value = "3.400.22"
$ws.Range("E3").Value = "  -4.37%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "578.09"
$ws.Range("E5").Value = "  -5.08%  "
$ws.Range("D6").Value = "131.57"
$ws.Range("E6").Value = "  -9.20%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.399.47"
$ws.Range("E8").Value = "  -4.48%  "
$ws.Range("E9").Value = "  -7.81%  "
$ws.Range("E10").Value = "  -10.48%  "
$ws.Range("D11").Value = "6.96"
$ws.Range("E11").Value = "  -10.89%  "
$ws.Range("E12").Value = "  -10.21%  "
$ws.Range("D13").Value = "3.976.07"
$ws.Range("E13").Value = "  -4.47%  "
$ws.Range("D14").Value = "0.0000176"
$ws.Range("E14").Value = "  -11.07%  "
$ws.Range("E15").Value = "  -2.01%  "
$ws.Range("D16").Value = "3.379.34"
$ws.Range("E16").Value = "  -4.91%  "
$ws.Range("E19").Value = "  -14.25%  "
$ws.Range("D20").Value = "5.64"
$ws.Range("E20").Value = "  -9.97%  "
$ws.Range("D21").Value = "13.45"
$ws.Range("E21").Value = "  -8.91%  "
$ws.Range("D22").Value = "376.91"
$ws.Range("E22").Value = "  -11.82%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("E24").Value = "  -10.36%  "
$ws.Range("D25").Value = "5.72"
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("D26").Value = "71.37"
$ws.Range("E26").Value = "  -8.34%  "
$ws.Range("D27").Value = "3.533.21"
$ws.Range("E27").Value = "  -4.47%  "
$ws.Range("E28").Value = "  -12.15%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").Value = "7.04"
$ws.Range("E30").Value = "  -12.52%  "
$ws.Range("E31").Value = "  -12.76%  "
$ws.Range("D32").Value = "7.93"
$ws.Range("E32").Value = "  -12.41%  "
$ws.Range("D33").Value = "3.414.65"
$ws.Range("E33").Value = "  -4.19%  "
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("D35").Value = "22.81"
$ws.Range("E35").Value = "  -6.98%  "
$ws.Range("E36").Value = "  -11.01%  "
$ws.Range("D37").Value = "171.03"
$ws.Range("E37").Value = "  -3.59%  "
$ws.Range("E38").Value = "  -14.16%  "
$ws.Range("D39").Value = "6.59"
$ws.Range("E39").Value = "  -13.98%  "
$ws.Range("E40").Value = "  -12.33%  "
$ws.Range("E41").Value = "  -14.43%  "
$ws.Range("D42").Value = "0.0753"
$ws.Range("E42").Value = "  -9.40%  "
$ws.Range("E43").Value = "  -8.57%  "
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").Value = "41.66"
$ws.Range("E45").Value = "  -8.66%  "
$ws.Range("E46").Value = "  -16.17%  "
$ws.Range("E47").Value = "  -11.70%  "
$ws.Range("E48").Value = "  -3.79%  "
$ws.Range("D49").Value = "21.95"
$ws.Range("E49").Value = "  -6.50%  "
$ws.Range("D50").Value = "6.45"
$ws.Range("E50").Value = "  -9.76%  "
$ws.Range("D51").Value = "2.179.61"
$ws.Range("E51").Value = "  -6.98%  "
# Row 17/18 swap: Avalanche moves to row 17, WrappedBTC moves to row 18
$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D17").Value = "25.82"
$ws.Range("E17").Value = "  -11.34%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "64.477.73"
$ws.Range("E18").Value = "  -3.34%  "
